# Applies the stat-bump + duplicate-row cleanup edit described by the diff.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (sheet1): simple F-column value bumps
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    2  = 801
    3  = 537
    4  = 284
    5  = 494
    6  = 1133
    7  = 319
    8  = 36
    10 = 114
    11 = 1150
    14 = 795
    15 = 818
    19 = 666
    20 = 193
    22 = 2364
    23 = 658
    24 = 69
    25 = 1912
    26 = 336
    27 = 2780
    28 = 513
    30 = 684
    32 = 102
    34 = 961
    35 = 1694
    36 = 336
    39 = 156
    41 = 157
    42 = 15
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# ---------------------------------------------------------------------------
# Sheet "演出" (sheet2): simple F-column value bumps
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 136
$ws2.Range("F11").Value = 16

# ---------------------------------------------------------------------------
# Sheet "全部类型" (sheet4): mirrors the same F-column bumps for the rows
# that correspond to "展览"/"演出" entries, PLUS removal of a duplicated
# row (row 18 duplicated row 17's "四月是你的谎言" concert) by deleting
# just the B:I cells of row 18 and shifting the rows below up.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    3  = 801
    4  = 537
    5  = 284
    6  = 494
    7  = 1133
    8  = 319
    9  = 36
    11 = 114
    12 = 1150
    14 = 795
    15 = 818
    17 = 136
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}

# Remove the duplicate row: delete B18:I18 and shift B19:I50 up by one,
# leaving column A (the static index column) untouched, then drop the
# now-empty trailing row.
$ws4.Range("B18:I18").Delete(-4162)  # -4162 = xlShiftUp

$wb.Save()
